$wb = $excel.ActiveWorkbook

$passengers = $wb.Worksheets.Item("Passengers")
$fuel = $wb.Worksheets.Item("Fuel")

# --- Passengers sheet: updated weighings (mass [kg]) ---
$passengers.Range("B2").Value = 80
$passengers.Range("B3").Value = 102
$passengers.Range("B4").Value = 82
$passengers.Range("B5").Value = 80
$passengers.Range("B6").Value = 85
$passengers.Range("B7").Value = 81
$passengers.Range("B8").Value = 65
$passengers.Range("B9").Value = 81
$passengers.Range("B10").Value = 100

# Unlock the input cells so the template stays fillable once sheet protection is on
$passengers.Range("B2:B10").Locked = $false

# --- Fuel sheet: updated block fuel figure ---
$fuel.Range("B2").Value = 2731
$fuel.Range("B2").Locked = $false

# Restore the cursor positions saved with the workbook (Fuel first, Passengers
# last so "Passengers" ends up as the active/visible tab, matching the source)
$fuel.Range("I12").Select() | Out-Null
$passengers.Range("K19").Select() | Out-Null

Write-Output "done"
